$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the rows that changed, per the commit
# "repull data, push all data, mean calculation".
$updates = @{
    3  = -4
    5  = 2
    6  = 4
    7  = 1
    9  = 1
    10 = 2
    11 = 1
    12 = -1
    13 = 8
    14 = -1
    15 = 1
    16 = 2
    17 = -1
    18 = 4
    19 = -2
    20 = 3
    21 = -8
    22 = -1
    23 = -1
    24 = -1
    26 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
